$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Checklist")

# The block of data rows 18-25 is being re-authored: row contents shift/merge
# and a brand-new row is inserted before the old "DUMMY_TestModuleCnt" /
# "ASDFClockTower" rows. Clear the whole old block first, then insert one row
# to make room, then (re)write every cell of the final state.

$ws.Range("A18:S25").ClearContents()
$ws.Rows("23").Insert()

# Row 18: TestData_07
$ws.Range("A18").Value = "TestData_07"
$ws.Range("B18").Value = 12345
$ws.Range("G18").Value = "X"
$ws.Range("H18").Value = "X"
$ws.Range("M18").Value = "ee_range"
$ws.Range("O18").Value = "description=- Component: Test`n- REPROG info: To be evaluated."

# Row 19: TestData_08
$ws.Range("A19").Value = "TestData_08"
$ws.Range("B19").Value = 12346
$ws.Range("I19").Value = "X"
$ws.Range("M19").Value = "ee_range"
$ws.Range("O19").Value = "description=- Component: TST Data`n- REPROG info: undefined"

# Row 20: TestData_09
$ws.Range("A20").Value = "TestData_09"
$ws.Range("B20").Value = 12347
$ws.Range("G20").Value = "X"
$ws.Range("M20").Value = "ee_range"
$ws.Range("O20").Value = "description=- Component: TST Data`n- REPROG info: tbd"

# Row 21: TestData_10
$ws.Range("A21").Value = "TestData_10"
$ws.Range("B21").Value = 12348
$ws.Range("H21").Value = "X"
$ws.Range("I21").Value = "X"
$ws.Range("M21").Value = "ee_range"
$ws.Range("O21").Value = "description=- Component: TST`n- REPROG info: t.b.d"

# Row 22: TestData_11
$ws.Range("A22").Value = "TestData_11"
$ws.Range("B22").Value = 12349
$ws.Range("H22").Value = "X"
$ws.Range("I22").Value = "X"
$ws.Range("M22").Value = "ee_range"
$ws.Range("O22").Value = "description=- Component: TST`n- REPROG info: use case REPROG must be set"

# Row 23: TestData_11 (new row, inserted)
$ws.Range("A23").Value = "TestData_11"
$ws.Range("B23").Value = 12349
$ws.Range("G23").Value = "X"
$ws.Range("M23").Value = "ee_range"
$ws.Range("O23").Value = "description=- Component: TST`n- REPROG info: use case REPROG must be set"

# Row 24: DUMMY_TestModuleCnt (shifted down from old row 23, content unchanged)
$ws.Range("A24").Value = "DUMMY_TestModuleCnt"
$ws.Range("B24").Value = 31416
$ws.Range("I24").Value = "X"
$ws.Range("M24").Value = "ee_erase"

# Row 25: ASDFClockTower (shifted down from old row 24, content unchanged)
$ws.Range("A25").Value = "ASDFClockTower"
$ws.Range("B25").Value = 111255
$ws.Range("G25").Value = "X"
$ws.Range("M25").Value = "ee_datablock"
$ws.Range("O25").Value = "description=- Component: ASDF`n- REPROG info: use case REPROG must not be set.`n- REPROG info: In certain cases there are two comments of this type.`nThere are also strings up to 160 characters per row, only on description fields and usually is not only one row. Like this example but a little bit longer."
